$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/1h change (E) columns for rows 2-42, 45-51
$ws.Range("D2").Value = "27.077.29"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "1.892.71"
$ws.Range("E3").Value = "  +1.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.77"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5175"
$ws.Range("E7").Value = "  +2.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3757"
$ws.Range("E8").Value = "  +2.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07215"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.12"
$ws.Range("E10").Value = "  +1.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8976"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07659"
$ws.Range("E12").Value = "  +1.74%  "

$ws.Range("D13").Value = "1.882.28"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.36"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.231"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008514"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.40"
$ws.Range("E18").Value = "  +1.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").Value = "27.128.94"
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.058"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").Value = "2.113.53"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  +1.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.417"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.287"
$ws.Range("E25").Value = "  +9.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.14"
$ws.Range("E26").Value = "  -1.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.733"
$ws.Range("E27").Value = "  -2.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.06"
$ws.Range("E28").Value = "  +0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.41"
$ws.Range("E29").Value = "  +0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.970"
$ws.Range("E30").Value = "  +6.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.786"
$ws.Range("E31").Value = "  +1.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09200"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05046"
$ws.Range("E33").Value = "  -2.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.238"
$ws.Range("E34").Value = "  +6.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7750"
$ws.Range("E35").Value = "  +2.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.978"
$ws.Range("E36").Value = "  -0.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.279"
$ws.Range("E37").Value = "  +1.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.595"
$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5602"
$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01989"
$ws.Range("E40").Value = "  -0.62%  "

$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.991"
$ws.Range("E42").Value = "  +5.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1514"
$ws.Range("E45").Value = "  +2.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4822"
$ws.Range("E46").Value = "  +1.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.19"
$ws.Range("E47").Value = "  +0.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("E49").Value = "  +1.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.39"
$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.84"
$ws.Range("E51").Value = "  +0.99%  "

# Row 43/44: Quant and FraxShare swap places (with updated price/volume)
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.632"
$ws.Range("E43").Value = "  +0.28%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.20"
$ws.Range("E44").Value = "  +2.98%  "
